$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23:49 down to 24:50
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new data
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = 44512
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 100112026
$ws.Range("G23").Value = "Haba"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 70
$ws.Range("K23").Value = 8000
$ws.Range("L23").Value = 8000
$ws.Range("M23").Value = 8000
$ws.Range("N23").Value = "$/saco 25 kilos"
$ws.Range("O23").Value = "Región del Maule"
$ws.Range("P23").Value = 320
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"

# Copy the date-number format style from the row above onto the new date cell
$ws.Range("D22").Copy()
$ws.Range("D23").PasteSpecial(-4122) # xlPasteFormats
